$d = $word.ActiveDocument
$d.Content.Find.Execute("<<court>>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<<courtName>>", 2)
